# Apply the "Deploy the implementation guide" update:
#  1. Metadata sheet: bump the generation Date value.
#  2. Concepts sheet: append a new concept row (SCID / Severe Combined
#     Immune Deficiency) at the bottom of the table, matching the
#     formatting of the preceding row.

$wb = $excel.ActiveWorkbook

# --- 1. Update the Date value on the Metadata sheet -----------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2023-11-09T20:21:55+00:00"

# --- 2. Append the new concept row on the Concepts sheet -------------------
$wsConcepts = $wb.Worksheets.Item("Concepts")

$lastRow = 17
$newRow = $lastRow + 1

# Copy the value of column A (the "1" level marker) as-is so it keeps its
# text type, then copy the whole row's formatting (borders/alignment/style)
# onto the new row without disturbing the shared workbook style table.
$wsConcepts.Range("A" + $lastRow).Copy() | Out-Null
$wsConcepts.Range("A" + $newRow).PasteSpecial(-4163) | Out-Null

$wsConcepts.Range("A" + $lastRow + ":D" + $lastRow).Copy() | Out-Null
$wsConcepts.Range("A" + $newRow + ":D" + $newRow).PasteSpecial(-4122) | Out-Null

$wsConcepts.Range("B" + $newRow).Value = "SCID"
$wsConcepts.Range("C" + $newRow).Value = "Severe Combined Immune Deficiency"
